# Improved the sample schedule:
#  - timezone changed from GMT to UTC
#  - uniqueid changed from leewkstest2 to poster-1
#  - selection / scroll position moved from M2:O2 (topLeftCell I1) to K2 (topLeftCell E1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "UTC"
$ws.Range("G2").Value = "poster-1"

$win = $excel.ActiveWindow
$win.ScrollColumn = 5
$win.ScrollRow = 1

$ws.Range("K2").Select()
